# mCSD 4.0.0 for TI
$wb = $excel.ActiveWorkbook

# --- Update Metadata sheet: Version and Date ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "4.0.0"
$meta.Range("B8").Value = "2025-05-21T13:41:24-05:00"

# --- Rename "Include #0" to "Include ValueSet #0" and replace its content ---
$inc0 = $wb.Worksheets.Item("Include #0")
$inc0.Name = "Include ValueSet #0"

# Clear old rows 3 and 4 (Codes / All codes / System URI rows no longer apply)
$inc0.Range("A3:B4").Clear()

# Set new content: ValueSet URL row
$inc0.Range("A1").Value = "ValueSet URL"
$inc0.Range("A2").Value = "https://profiles.ihe.net/ITI/mCSD/ValueSet/MCSDOrgDocSharingAffTypesVS"
